# Updates cryptos list (row 2-51, columns D price / E volume%) with refreshed
# market data. Rows 44/45 (Arweave/Cosmos) also swap places with new values,
# matching Cosmos's higher rank on this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='66.908.76'; E='  -3.66%  ' },
    @{ Row=3; D='3.720.24'; E='  +0.22%  ' },
    @{ Row=4; E='  -0.16%  ' },
    @{ Row=5; D='588.30'; E='  -3.72%  ' },
    @{ Row=6; D='171.28'; E='  -3.72%  ' },
    @{ Row=7; D='3.720.33'; E='  +0.26%  ' },
    @{ Row=8; E='  -0.07%  ' },
    @{ Row=9; D='0.517'; E='  -1.90%  ' },
    @{ Row=10; D='0.157'; E='  -4.56%  ' },
    @{ Row=11; D='6.24'; E='  -4.58%  ' },
    @{ Row=12; D='0.456'; E='  -4.52%  ' },
    @{ Row=13; D='37.43'; E='  -5.40%  ' },
    @{ Row=14; D='0.0000241'; E='  -4.48%  ' },
    @{ Row=15; D='4.349.89'; E='  +0.42%  ' },
    @{ Row=16; D='3.726.65'; E='  +0.29%  ' },
    @{ Row=17; D='67.098.19'; E='  -3.49%  ' },
    @{ Row=18; E='  -4.64%  ' },
    @{ Row=19; D='7.03'; E='  -5.65%  ' },
    @{ Row=20; D='15.92'; E='  -1.94%  ' },
    @{ Row=21; D='482.39'; E='  -3.41%  ' },
    @{ Row=22; D='9.01'; E='  -1.03%  ' },
    @{ Row=23; D='0.716'; E='  +0.05%  ' },
    @{ Row=24; D='83.43'; E='  -2.80%  ' },
    @{ Row=25; D='2.34'; E='  -8.98%  ' },
    @{ Row=26; D='0.0000134'; E='  +1.07%  ' },
    @{ Row=27; D='12.11'; E='  -5.77%  ' },
    @{ Row=28; D='10.03'; E='  -10.27%  ' },
    @{ Row=29; D='1.00'; E='  +0.41%  ' },
    @{ Row=30; D='2.88'; E='  -0.64%  ' },
    @{ Row=31; D='2.37'; E='  -3.18%  ' },
    @{ Row=32; D='31.77'; E='  +5.12%  ' },
    @{ Row=33; D='7.63'; E='  -4.55%  ' },
    @{ Row=34; D='0.107'; E='  -5.25%  ' },
    @{ Row=35; E='  +0.07%  ' },
    @{ Row=36; E='  -4.42%  ' },
    @{ Row=37; D='0.134'; E='  -2.88%  ' },
    @{ Row=38; D='5.65'; E='  -6.96%  ' },
    @{ Row=39; D='0.320'; E='  -7.86%  ' },
    @{ Row=40; D='442.56'; E='  +2.51%  ' },
    @{ Row=41; D='48.67'; E='  -1.99%  ' },
    @{ Row=42; D='1.95'; E='  -4.85%  ' },
    @{ Row=43; D='2.83'; E='  -6.10%  ' },
    @{ Row=44; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.15'; E='  -4.33%  ' },
    @{ Row=45; B='Arweave'; C='https://coinranking.com/coin/7XWg41D1+arweave-ar'; D='40.94'; E='  -9.94%  ' },
    @{ Row=46; D='2.783.05'; E='  -5.24%  ' },
    @{ Row=47; E='  +0.05%  ' },
    @{ Row=48; D='138.55'; E='  +0.29%  ' },
    @{ Row=49; D='0.0343'; E='  -4.35%  ' },
    @{ Row=50; D='25.57'; E='  -5.26%  ' },
    @{ Row=51; D='22.75'; E='  +7.44%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$row").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$row").Value = $u.C }
    if ($u.ContainsKey('D')) {
        # Prefix with an apostrophe so Excel stores the price as text (matching
        # the original inline-string cells) instead of inferring a number and
        # mangling values like "588.30" -> 588.29999999999995 or dropping
        # trailing zeros / thousands separators.
        $ws.Range("D$row").Value = "'" + $u.D
        # Clear the resulting quote-prefix style so the cell's applied style
        # stays the same as before (no style index) rather than picking up a
        # "number stored as text" marker.
        $ws.Range("D$row").Style = "Normal"
    }
    if ($u.ContainsKey('E')) { $ws.Range("E$row").Value = $u.E }
}
